$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, report week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Cells that change between numeric and text representation (need style carried over) ---
$ws.Range("D14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D14").PasteSpecial(-4122)

$ws.Range("E14").Value = "'***.*"
$ws.Range("F14").Copy()
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("C16").Value = 2
$ws.Range("F16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 2

$ws.Range("D16").Value = 1
$ws.Range("G16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1

$ws.Range("E16").Value = 100
$ws.Range("H16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = 100

$ws.Range("D17").Value = 1
$ws.Range("C17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 1

$ws.Range("E17").Value = 0
$ws.Range("H17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = 0

$ws.Range("C18").Value = "'0"
$ws.Range("A18").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D18").Value = "'0"
$ws.Range("A18").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").Value = "'***.*"
$ws.Range("A18").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("D28").Value = "'0"
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = "'***.*"
$ws.Range("F28").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("D29").Value = "'0"
$ws.Range("C29").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").Value = "'***.*"
$ws.Range("F29").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Plain numeric value updates ---
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 22
$ws.Range("J16").Value = 11
$ws.Range("L16").Value = 69.230769230769
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -65.625
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 150
$ws.Range("I17").Value = 51
$ws.Range("J17").Value = 32
$ws.Range("K17").Value = 59.375
$ws.Range("L17").Value = -19.047619047619
$ws.Range("M17").Value = 4.081632653061
$ws.Range("N17").Value = -53.636363636363
$ws.Range("F18").Value = 3
$ws.Range("H18").Value = 200
$ws.Range("M18").Value = -60.747663551401
$ws.Range("N18").Value = -87.537091988130
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 166.666666666667
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -29.411764705882
$ws.Range("I19").Value = 270
$ws.Range("J19").Value = 198
$ws.Range("K19").Value = 36.363636363636
$ws.Range("L19").Value = 50
$ws.Range("M19").Value = 84.931506849315
$ws.Range("N19").Value = 12.033195020746
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("I20").Value = 120
$ws.Range("J20").Value = 61
$ws.Range("K20").Value = 96.721311475409
$ws.Range("L20").Value = 166.666666666667
$ws.Range("M20").Value = 215.789473684211
$ws.Range("N20").Value = -83.216783216783
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = 71.428571428571
$ws.Range("F21").Value = 44
$ws.Range("G21").Value = 49
$ws.Range("H21").Value = -10.204081632653
$ws.Range("I21").Value = 509
$ws.Range("J21").Value = 329
$ws.Range("K21").Value = 54.711246200607
$ws.Range("L21").Value = 46.264367816092
$ws.Range("M21").Value = 39.071038251366
$ws.Range("N21").Value = -65.421195652173
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -41.666666666666
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 37
$ws.Range("H24").Value = 10.810810810810
$ws.Range("I24").Value = 488
$ws.Range("J24").Value = 287
$ws.Range("K24").Value = 70.034843205574
$ws.Range("L24").Value = 71.830985915493
$ws.Range("M24").Value = -10.294117647058
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = -9.090909090909
$ws.Range("I25").Value = 183
$ws.Range("J25").Value = 142
$ws.Range("K25").Value = 28.873239436619
$ws.Range("L25").Value = 36.567164179104
$ws.Range("M25").Value = -15.668202764977
